$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.554.49"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.806.01"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "2.060.30"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "1.800.66"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.644"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "34.577.69"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "251.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "0.0₃0799"
$ws.Range("E20").Value = "  +7.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.45%  "
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0526"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("B33").Value = "Swop.fi"
$ws.Range("C33").Value = "https://coinranking.com/coin/yrCr2HW2c+swopfi-swop"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "504.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +867.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").Value = "1.426.10"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.644"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.966"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "82.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("D49").Value = "1.955.87"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.22%  "
$ws.Range("E51").Value = "  -0.11%  "
